$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text
$ws.Range("D1").Value = "Memory Usage (bytes)"

# Update Run Time (ms) values in column C
$ws.Range("C2").Value = 20.9958553314209
$ws.Range("C3").Value = 17.98081398010254
$ws.Range("C4").Value = 17.68088340759277
$ws.Range("C5").Value = 17.72093772888184
$ws.Range("C6").Value = 18.0511474609375
